$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue "D2" "307.69"
Set-TextValue "E2" "3.19%"
Set-TextValue "D3" "36.27"
Set-TextValue "E3" "3.38%"
Set-TextValue "D4" "5.149"
Set-TextValue "E4" "3.70%"
Set-TextValue "D5" "0.08143"
Set-TextValue "E5" "4.53%"
Set-TextValue "D6" "1.931"
Set-TextValue "E6" "2.70%"
Set-TextValue "D7" "7.776"
Set-TextValue "E7" "0.81%"
Set-TextValue "D8" "0.9318"
Set-TextValue "E8" "1.12%"
Set-TextValue "D9" "0.1378"
Set-TextValue "E9" "28.42%"
Set-TextValue "D10" "0.1932"
Set-TextValue "E10" "6.58%"
Set-TextValue "D11" "0.09265"
Set-TextValue "E11" "-0.02%"
Set-TextValue "D12" "0.03559"
Set-TextValue "E12" "0.46%"
Set-TextValue "D13" "0.09869"
Set-TextValue "E13" "-0.10%"
Set-TextValue "D14" "0.001406"
Set-TextValue "E14" "0.73%"
Set-TextValue "D15" "0.005905"
Set-TextValue "E15" "3.64%"
Set-TextValue "D16" "3.553"
Set-TextValue "E16" "2.03%"
Set-TextValue "D17" "4.184"
Set-TextValue "E17" "4.10%"
Set-TextValue "D18" "2.970"
Set-TextValue "E18" "2.09%"
Set-TextValue "D19" "0.3441"
Set-TextValue "E19" "-0.02%"
Set-TextValue "D20" "0.1305"
Set-TextValue "E20" "1.15%"
Set-TextValue "D21" "4.902"
Set-TextValue "E21" "-2.72%"
Set-TextValue "D22" "0.2498"
Set-TextValue "E22" "7.89%"
Set-TextValue "D23" "0.04548"
Set-TextValue "E23" "0.03%"
Set-TextValue "D24" "0.001214"
Set-TextValue "E24" "0.06%"
Set-TextValue "D25" "0.004887"
Set-TextValue "E25" "6.53%"
Set-TextValue "D26" "0.0001241"
Set-TextValue "E26" "-0.79%"
Set-TextValue "E27" "5.94%"
Set-TextValue "D39" "0.02001"
Set-TextValue "E39" "7.03%"
Set-TextValue "D40" "0.04932"
Set-TextValue "E40" "5.65%"
Set-TextValue "D41" "0.01112"
Set-TextValue "E41" "15.99%"
Set-TextValue "D42" "0.007665"
Set-TextValue "E42" "0.91%"
Set-TextValue "D43" "0.1379"
Set-TextValue "E43" "4.50%"
Set-TextValue "D44" "0.002101"
Set-TextValue "E44" "-0.94%"
Set-TextValue "D45" "0.01053"
Set-TextValue "E45" "-5.68%"
Set-TextValue "D46" "0.00006465"
Set-TextValue "E46" "7.57%"
Set-TextValue "D47" "0.00000000751"
Set-TextValue "E47" "0.01%"
Set-TextValue "E48" "0.33%"
Set-TextValue "E49" "-8.68%"
Set-TextValue "D50" "0.00002101"
Set-TextValue "E50" "0.01%"
Set-TextValue "D51" "0.0002001"
Set-TextValue "E51" "0.01%"
